$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 75.17253906259769
$ws.Range("H2").Value = 4.188286731620309

$ws.Range("G3").Value = 41.45877320543772
$ws.Range("H3").Value = 1.995496902034478

$ws.Range("G4").Value = 32.82951004874731
$ws.Range("H4").Value = 1.477418817867068

$ws.Range("G5").Value = 38.99693567293362
$ws.Range("H5").Value = 1.786366103883976

$ws.Range("G6").Value = 42.62121031591623
$ws.Range("H6").Value = 2.741089328914315

$ws.Range("G7").Value = 90.17548792134531
$ws.Range("H7").Value = 0.8001611251767615

$ws.Range("G8").Value = 45.79194179155397
$ws.Range("H8").Value = 4.391401271657076

$ws.Range("G9").Value = 42.4338049218881
$ws.Range("H9").Value = 4.421035565841859

$ws.Range("G10").Value = 30.85314796492898
$ws.Range("H10").Value = 2.753833965186642

$ws.Range("G11").Value = 50.17943710733019
$ws.Range("H11").Value = 3.049211438980877

$ws.Range("G12").Value = 42.92282107595251
$ws.Range("H12").Value = 3.359724299501527

$ws.Range("G13").Value = 33.68210803386253
$ws.Range("H13").Value = 4.184512419492172
